# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" quarter sheet (cloned from the existing "2022-Q2"
# sheet so it keeps the same layout/styling) right after "总计", fills it
# with the new quarter's figures, and prepends a matching row to the "总计"
# summary sheet. The other existing quarter sheets (2022-Q2, 2022-Q1,
# 2021-Q3, 2021-Q2) are left untouched - they simply shift one tab to the
# right to make room for the new one.

$wb = $excel.ActiveWorkbook

# --- 1. Clone the "2022-Q2" sheet to use as the template for "2022-Q3" ---
$template = $wb.Worksheets.Item(2)
$template.Copy($template)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# --- 2. Overwrite the new sheet's data row with the 2022-Q3 figures ---
# (headers, A2, B2 fund code and C2 fund name are already correct from the
# clone). D2:G2 are stored as text in this workbook, so force text with a
# leading apostrophe (otherwise Excel would coerce the numeric-looking
# string into a number).
$q3.Range("D2").Value = "'27.03"
$q3.Range("E2").Value = "'99.07"
$q3.Range("F2").Value = "'4.53"
$q3.Range("G2").Value = "'1.2245"
$q3.Range("H2").Value = 6

# --- 3. Restore the originally-selected tab (the last sheet, "2021-Q2") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

# --- 4. Prepend a 2022-Q3 row to the "总计" summary sheet, pushing the ---
# --- existing rows down by one.                                      ---
$summary = $wb.Worksheets.Item(1)

# Grab the formatted/styled template cell (A2) before we touch anything,
# then stamp it onto the brand-new row 6 so it picks up the same style as
# the existing index column cells.
$summary.Range("A2").Copy($summary.Range("A6"))

$summary.Range("B6").Value = "2021-Q2"
$summary.Range("C6").Value = 1
$summary.Range("D6").Value = 1.46

$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 1.14

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 1.24

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 1.2

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 1.22

$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
